# Add PSP Time Recording Log entries for rows 22-25 (new activities logged
# by the author after 10/17), matching the author's "Add file via uploading"
# re-upload of the tracked workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: 10/22 - survey data organizing
$ws.Range("A22").Value = "10월 22일"
$ws.Range("B22").Value = 0.77083333333333337
$ws.Range("C22").Value = 0.875
$ws.Range("D22").Value = 20
$ws.Range("E22").Value = 150
$ws.Range("F22").Value = "설문지 자료 정리"

# Row 23: 10/24 - build initial data set from survey data
$ws.Range("A23").Value = "10월 24일"
$ws.Range("B23").Value = 0.45833333333333331
$ws.Range("C23").Value = 0.78402777777777777
$ws.Range("D23").Value = 70
$ws.Range("E23").Value = 469
$ws.Range("F23").Value = "설문지 자료 토대로 Initial Data set 작성"

# Row 24: 10/29 - web lecture 4
$ws.Range("A24").Value = "10월 29일"
$ws.Range("B24").Value = 0.70833333333333337
$ws.Range("C24").Value = 0.79166666666666663
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 120
$ws.Range("F24").Value = "web1 4강"

# Row 25: 11/3 - web lecture 6
$ws.Range("A25").Value = "11월 3일"
$ws.Range("B25").Value = 0.5
$ws.Range("C25").Value = 0.16388888888888889
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 236
$ws.Range("F25").Value = "web1 6강"

# Move the active selection to F26, matching where the author left off editing
$ws.Range("F26").Select()
